$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp (row 1 title) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 13:46"

# --- Refresh country statistics (cols B,C,D,E,F,G,H) for the rows whose
#     numbers moved between the two data pulls ---
$ws.Range("B4").Value = 3480059
$ws.Range("C4").Value = 576
$ws.Range("D4").Value = 1550324
$ws.Range("E4").Value = 1791462
$ws.Range("G4").Value = 26
$ws.Range("H4").Value = 138273

$ws.Range("B6").Value = 911629
$ws.Range("C6").Value = 3984
$ws.Range("D6").Value = 573953
$ws.Range("E6").Value = 313888
$ws.Range("G6").Value = 61
$ws.Range("H6").Value = 23788

$ws.Range("B32").Value = 65269
$ws.Range("C32").Value = 155
$ws.Range("D32").Value = 55799
$ws.Range("E32").Value = 8996
$ws.Range("G32").Value = 6
$ws.Range("H32").Value = 474

$ws.Range("B52").Value = 33016
$ws.Range("C52").Value = 70
$ws.Range("E52").Value = 1448

$ws.Range("B62").Value = 19021
$ws.Range("C62").Value = 73
$ws.Range("D62").Value = 17073
$ws.Range("E62").Value = 1239
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 709

$ws.Range("B64").Value = 17061
$ws.Range("C64").Value = 116
$ws.Range("D64").Value = 10328
$ws.Range("E64").Value = 6695

$ws.Range("D67").Value = 8074
$ws.Range("E67").Value = 5733
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 65

$ws.Range("B95").Value = 5343
$ws.Range("C95").Value = 263
$ws.Range("D95").Value = 2646
$ws.Range("E95").Value = 2658
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = 39

$ws.Range("B110").Value = 2651
$ws.Range("C110").Value = 5
$ws.Range("E110").Value = 652

$ws.Range("B119").Value = 1905
$ws.Range("C119").Value = 5
$ws.Range("D119").Value = 1882
$ws.Range("E119").Value = 13

# --- Re-sort "Lesoto" / "Martinica" alphabetically and swap in the
#     (refreshed) stats that go with each name ---
$ws.Range("A169").Value = "Lesoto"
$ws.Range("B169").Value = 256
$ws.Range("C169").Value = 11
$ws.Range("D169").Value = 48
$ws.Range("E169").Value = 205
$ws.Range("H169").Value = 3

$ws.Range("A170").Value = "Martinica"
$ws.Range("B170").Value = 255
$ws.Range("C170").Value = 0
$ws.Range("D170").Value = 98
$ws.Range("E170").Value = 142
$ws.Range("H170").Value = 15
